# Add a new row of data for Triton to the Table_of_Satellites workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: Satellite name in column A, Rc(km) in C, Rc/Rp in D, Mc/Mp in E.
# Columns B, F, G and H are left untouched/blank, same as the source diff.
# Only the four needed cells are created, each copying the formatting of the
# corresponding cell one row above, so the existing style/font entries are
# reused instead of new ones being created.
foreach ($col in @("A", "C", "D", "E")) {
    $srcAddr = $col + "18"
    $dstAddr = $col + "19"
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

$ws.Range("A19").Value = "Triton"
$ws.Range("C19").Value = 1028
$ws.Range("D19").Value = 0.76
$ws.Range("E19").Value = 0.75

# Match the row height (18pt) used by all the other data rows in the table.
$ws.Rows.Item(19).RowHeight = 18

# Move the active selection to E19, matching the saved cursor position.
$ws.Range("E19").Select()
